$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33389.25
$ws.Range("J3").Value = 33389.25
$ws.Range("L3").Value = 33389.25
$ws.Range("N3").Value = -33617.25
$ws.Range("H15").Value = 1685844.9
$ws.Range("I15").Value = 1685844.9
$ws.Range("K15").Value = 5057534.699999999
$ws.Range("M15").Value = -5057365.699999999
$ws.Range("H102").Value = 33389.25
$ws.Range("J102").Value = 33389.25
$ws.Range("L102").Value = 33389.25
$ws.Range("N102").Value = -39879.25
$ws.Range("H132").Value = 3059.842
$ws.Range("I132").Value = 3165.1538
$ws.Range("J132").Value = 2831.6667
$ws.Range("K132").Value = 9495.4614
$ws.Range("L132").Value = 8495.000100000001
$ws.Range("M132").Value = -6965.4614
$ws.Range("N132").Value = -13555.0001
$ws.Range("H137").Value = 20006.79
$ws.Range("I137").Value = 43865.59
$ws.Range("J137").Value = 8418.228999999999
$ws.Range("K137").Value = 131596.77
$ws.Range("L137").Value = 25254.687
$ws.Range("M137").Value = -129046.77
$ws.Range("N137").Value = -30354.687
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7701771
$ws.Range("I32").Value = 8071138.5
$ws.Range("K32").Value = 8071138.5
$ws.Range("M32").Value = -8070851.5
$ws.Range("H88").Value = 1747.3334
$ws.Range("I88").Value = 1393.6666
$ws.Range("J88").Value = 2101
$ws.Range("K88").Value = 1393.6666
$ws.Range("L88").Value = 2101
$ws.Range("M88").Value = -987.6666
$ws.Range("N88").Value = -2913
$ws.Range("H91").Value = 1747.3334
$ws.Range("I91").Value = 1393.6666
$ws.Range("J91").Value = 2101
$ws.Range("K91").Value = 1393.6666
$ws.Range("L91").Value = 2101
$ws.Range("M91").Value = 10.33339999999998
$ws.Range("N91").Value = -4909
$ws.Range("H122").Value = 2317.861
$ws.Range("I122").Value = 1278.75
$ws.Range("K122").Value = 3836.25
$ws.Range("M122").Value = -1386.25
$ws.Range("H123").Value = 120429
$ws.Range("J123").Value = 120429
$ws.Range("L123").Value = 120429
$ws.Range("N123").Value = -130229
$ws.Range("H132").Value = 6698.7144
$ws.Range("I132").Value = 2357.889
$ws.Range("J132").Value = 14512.2
$ws.Range("K132").Value = 7073.667
$ws.Range("L132").Value = 43536.60000000001
$ws.Range("M132").Value = -4543.667
$ws.Range("N132").Value = -48596.60000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 48260.8
$ws.Range("J81").Value = 48260.8
$ws.Range("L81").Value = 48260.8
$ws.Range("N81").Value = -50382.8
$ws.Range("H84").Value = 48260.8
$ws.Range("J84").Value = 48260.8
$ws.Range("L84").Value = 144782.4
$ws.Range("N84").Value = -155390.4
$ws.Range("H100").Value = 29397.834
$ws.Range("J100").Value = 29397.834
$ws.Range("L100").Value = 29397.834
$ws.Range("N100").Value = -31561.834
$ws.Range("H107").Value = 1277.75
$ws.Range("I107").Value = 1289.6666
$ws.Range("K107").Value = 1289.6666
$ws.Range("M107").Value = 630.3334
$ws.Range("H110").Value = 317825.25
$ws.Range("J110").Value = 317825.25
$ws.Range("L110").Value = 317825.25
$ws.Range("N110").Value = -326005.25
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 524004.44
$ws.Range("I31").Value = 11230.823
$ws.Range("J31").Value = 835331.25
$ws.Range("K31").Value = 11230.823
$ws.Range("L31").Value = 835331.25
$ws.Range("M31").Value = -10935.823
$ws.Range("N31").Value = -835921.25
$ws.Range("H34").Value = 524004.44
$ws.Range("I34").Value = 11230.823
$ws.Range("J34").Value = 835331.25
$ws.Range("K34").Value = 11230.823
$ws.Range("L34").Value = 835331.25
$ws.Range("M34").Value = -11028.823
$ws.Range("N34").Value = -835735.25
$ws.Range("H107").Value = 717.7857
$ws.Range("I107").Value = 679.9
$ws.Range("K107").Value = 679.9
$ws.Range("M107").Value = 1240.1
$ws.Range("H134").Value = 916719
$ws.Range("I134").Value = 1117988.4
$ws.Range("K134").Value = 3353965.2
$ws.Range("M134").Value = -3351430.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12675.906
$ws.Range("I2").Value = 195.77272
$ws.Range("K2").Value = 1174.63632
$ws.Range("M2").Value = -1061.63632
$ws.Range("H17").Value = 3666.3333
$ws.Range("J17").Value = 3666.3333
$ws.Range("L17").Value = 10998.9999
$ws.Range("N17").Value = -11336.9999
$ws.Range("H86").Value = 2457.6667
$ws.Range("I86").Value = 898.4
$ws.Range("K86").Value = 2695.2
$ws.Range("M86").Value = -1509.2
$ws.Range("H89").Value = 2457.6667
$ws.Range("I89").Value = 898.4
$ws.Range("K89").Value = 8085.599999999999
$ws.Range("M89").Value = -2157.599999999999
$ws.Range("H104").Value = 3845.7
$ws.Range("I104").Value = 3892.8
$ws.Range("J104").Value = 3798.6
$ws.Range("K104").Value = 11678.4
$ws.Range("L104").Value = 11395.8
$ws.Range("M104").Value = -9057.400000000001
$ws.Range("N104").Value = -16637.8
$ws.Range("H109").Value = 1983.5714
$ws.Range("I109").Value = 1283.3334
$ws.Range("K109").Value = 3850.0002
$ws.Range("M109").Value = -2810.0002
$ws.Range("H111").Value = 943.5
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 885.2
$ws.Range("I114").Value = 812.5
$ws.Range("J114").Value = 911.63635
$ws.Range("K114").Value = 2437.5
$ws.Range("L114").Value = 2734.90905
$ws.Range("M114").Value = 816.5
$ws.Range("N114").Value = -9242.90905
$ws.Range("H137").Value = 5821.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 188597.12
$ws.Range("I140").Value = 188597.12
$ws.Range("K140").Value = 565791.36
$ws.Range("M140").Value = -560611.36
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19999
$ws.Range("J46").Value = 19999
$ws.Range("L46").Value = 19999
$ws.Range("N46").Value = -20311
$ws.Range("H57").Value = 60000
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 100000
$ws.Range("N57").Value = -101640
$ws.Range("H102").Value = 2934.2856
$ws.Range("I102").Value = 2360.5334
$ws.Range("K102").Value = 2360.5334
$ws.Range("M102").Value = -738.5333999999998
$ws.Range("H107").Value = 871.5714
$ws.Range("I107").Value = 733.5
$ws.Range("J107").Value = 1700
$ws.Range("K107").Value = 733.5
$ws.Range("L107").Value = 1700
$ws.Range("M107").Value = 1186.5
$ws.Range("N107").Value = -5540
$ws.Range("H122").Value = 6777.7144
$ws.Range("I122").Value = 5391.263
$ws.Range("K122").Value = 16173.789
$ws.Range("M122").Value = -13723.789
$ws.Range("H123").Value = 51507.75
$ws.Range("J123").Value = 51507.75
$ws.Range("L123").Value = 51507.75
$ws.Range("N123").Value = -56407.75
$ws.Range("H124").Value = 96779.664
$ws.Range("J124").Value = 96779.664
$ws.Range("L124").Value = 96779.664
$ws.Range("N124").Value = -106599.664
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4201
$ws.Range("I61").Value = 4201
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4201
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3999
$ws.Range("N61").ClearContents()
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H113").Value = 4201
$ws.Range("I113").Value = 4201
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4201
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2031
$ws.Range("N113").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H107").Value = 62501624
$ws.Range("J107").Value = 1001.5
$ws.Range("L107").Value = 3004.5
$ws.Range("N107").Value = -6844.5
$ws.Range("H136").Value = 4577.2856
$ws.Range("I136").Value = 3908.3
$ws.Range("K136").Value = 11724.9
$ws.Range("M136").Value = -9174.900000000001
